$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row (data occupies rows 1..94 in columns A:D,
# but compute dynamically in case it ever differs).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $cVal = $cCell.Value2
    $dVal = $dCell.Value2
    $cCell.Value2 = $dVal
    $dCell.Value2 = $cVal
}
